# Corrección algoritmo de análisis
#
# Appends a closing period "." as its own run right after the run that
# ends in "... más altos de la variable dependiente, (siendo esta la
# concentración)" — matching the target OOXML, where the new "." gets its
# own <w:r><w:rPr> (Segoe UI, sz/szCs 24) rather than being folded into
# the previous run's text node.

$d = $word.ActiveDocument

$anchorText = "más altos de la variable dependiente, (siendo esta la concentración)"

$find = $d.Content.Find
$find.ClearFormatting()
[void]$find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($find.Found) {
    $targetRng = $find.Parent
    $paraRng = $targetRng.Paragraphs(1).Range

    # Grab a "." run already formatted exactly like the surrounding text
    # (Segoe UI, sz 24 / szCs 24) from earlier in the very same paragraph
    # ("... código de departamento. Ya definidas ...") so the inserted
    # period ends up with the identical <w:rPr> (including w:cs), instead
    # of the bare run InsertAfter()/TypeText() would otherwise produce.
    $srcFind = $paraRng.Duplicate.Find
    $srcFind.ClearFormatting()
    [void]$srcFind.Execute("departamento.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

    $periodSrc = $srcFind.Parent.Duplicate
    $periodSrc.Start = $srcFind.Parent.End - 1
    $periodSrc.End = $srcFind.Parent.End
    $formattedPeriod = $periodSrc.FormattedText

    # Collapse to the end of the matched text (right before the paragraph
    # mark) and insert a placeholder run there.
    $targetRng.Collapse(0)
    $insertStart = $targetRng.Start
    $targetRng.InsertAfter(".")

    # Re-stamp that new run with the copied, fully-formatted "." so it
    # carries the same rFonts/sz/szCs as its neighbour.
    $newRunRng = $d.Range($insertStart, $insertStart + 1)
    $newRunRng.FormattedText = $formattedPeriod
}
